$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "data"
$dataSheet = $wb.Worksheets.Item("Sheet1")
$dataSheet.Name = "data"

# Add new "environment" sheet after "data"
$envSheet = $wb.Worksheets.Add()
$envSheet.Name = "environment"
$envSheet.Range("A1").Value = "environment 124.7"
$envSheet.Move($null, $wb.Worksheets.Item("data"))

# Update header row on "data" sheet
$dataSheet.Range("B1").Value = "luminance1"
$dataSheet.Range("C1").Value = "luminance2"
$dataSheet.Range("D1").Value = "luminance3"
$dataSheet.Range("E1").Value = "luminance_average"

# Restore "data" as the active sheet/selection (scrolled back to top, cell G5 selected)
$dataSheet.Activate()
[void]$dataSheet.Range("G5").Select()
